$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1053.25
$ws.Range("I41").Value = 249
$ws.Range("K41").Value = 249
$ws.Range("M41").Value = 191

$ws.Range("H62").Value = 8385.549999999999
$ws.Range("I62").Value = 7145.7856
$ws.Range("J62").Value = 11278.333
$ws.Range("K62").Value = 7145.7856
$ws.Range("L62").Value = 11278.333
$ws.Range("M62").Value = -6521.7856
$ws.Range("N62").Value = -12526.333

$ws.Range("H65").Value = 8385.549999999999
$ws.Range("I65").Value = 7145.7856
$ws.Range("J65").Value = 11278.333
$ws.Range("K65").Value = 35728.928
$ws.Range("L65").Value = 56391.665
$ws.Range("M65").Value = -32608.928
$ws.Range("N65").Value = -62631.665

$ws.Range("H74").Value = 9546.308000000001
$ws.Range("I74").Value = 9410.299999999999
$ws.Range("K74").Value = 9410.299999999999
$ws.Range("M74").Value = -8474.299999999999

$ws.Range("H77").Value = 9546.308000000001
$ws.Range("I77").Value = 9410.299999999999
$ws.Range("K77").Value = 47051.5
$ws.Range("M77").Value = -42371.5

$ws.Range("H97").Value = 1476.25
$ws.Range("J97").Value = 1476.25
$ws.Range("L97").Value = 4428.75
$ws.Range("N97").Value = -5420.75

$ws.Range("H137").Value = 16187.714
$ws.Range("I137").Value = 31710.9
$ws.Range("J137").Value = 2075.7273
$ws.Range("K137").Value = 95132.70000000001
$ws.Range("L137").Value = 6227.1819
$ws.Range("M137").Value = -92582.70000000001
$ws.Range("N137").Value = -11327.1819

$ws.Range("H138").Value = 25570.373
$ws.Range("J138").Value = 58199.723
$ws.Range("L138").Value = 174599.169
$ws.Range("N138").Value = -184879.169

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22763.51
$ws.Range("I32").Value = 23686.377
$ws.Range("K32").Value = 23686.377
$ws.Range("M32").Value = -23399.377

$ws.Range("H74").Value = 279457.88
$ws.Range("I74").Value = 601422
$ws.Range("J74").Value = 11154.417
$ws.Range("K74").Value = 601422
$ws.Range("L74").Value = 11154.417
$ws.Range("M74").Value = -600548
$ws.Range("N74").Value = -12902.417

$ws.Range("H77").Value = 279457.88
$ws.Range("I77").Value = 601422
$ws.Range("J77").Value = 11154.417
$ws.Range("K77").Value = 3007110
$ws.Range("L77").Value = 55772.085
$ws.Range("M77").Value = -3002742
$ws.Range("N77").Value = -64508.085

$ws.Range("H132").Value = 1989.4546
$ws.Range("I132").Value = 1998.4
$ws.Range("J132").Value = 1900
$ws.Range("K132").Value = 5995.200000000001
$ws.Range("L132").Value = 5700
$ws.Range("M132").Value = -3465.200000000001
$ws.Range("N132").Value = -10760

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2030.3
$ws.Range("I105").Value = 1459.6666
$ws.Range("K105").Value = 1459.6666
$ws.Range("M105").Value = 287.3334

$ws.Range("H107").Value = 3572.5757
$ws.Range("I107").Value = 3623.4546
$ws.Range("J107").Value = 3470.818
$ws.Range("K107").Value = 3623.4546
$ws.Range("L107").Value = 3470.818
$ws.Range("M107").Value = -1703.4546
$ws.Range("N107").Value = -7310.818

$ws.Range("H134").Value = 2771.4443
$ws.Range("I134").Value = 1414
$ws.Range("K134").Value = 4242
$ws.Range("M134").Value = -1707

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 8690.538
$ws.Range("I99").Value = 5748
$ws.Range("K99").Value = 5748
$ws.Range("M99").Value = -4250

$ws.Range("H126").Value = 8690.538
$ws.Range("I126").Value = 5748
$ws.Range("K126").Value = 17244
$ws.Range("M126").Value = -14774

$ws.Range("H132").Value = 41361.12
$ws.Range("I132").Value = 46364.953
$ws.Range("K132").Value = 139094.859
$ws.Range("M132").Value = -136564.859

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 3000
$ws.Range("J48").Value = 3000
$ws.Range("L48").Value = 9000
$ws.Range("N48").Value = -9500

$ws.Range("H63").Value = 5000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 15000
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -16498

$ws.Range("H66").Value = 5000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 45000
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -52488

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 37999
$ws.Range("J101").Value = 37999
$ws.Range("L101").Value = 37999
$ws.Range("N101").Value = -44489

$ws.Range("H105").Value = 5300335.5
$ws.Range("J105").Value = 5300335.5
$ws.Range("L105").Value = 5300335.5
$ws.Range("N105").Value = -5307323.5

$ws.Range("H107").Value = 494.83334
$ws.Range("J107").Value = 699.625
$ws.Range("L107").Value = 699.625
$ws.Range("N107").Value = -4539.625

$ws.Range("H132").Value = 2536.4075
$ws.Range("I132").Value = 2200.4707
$ws.Range("J132").Value = 3107.5
$ws.Range("K132").Value = 6601.4121
$ws.Range("L132").Value = 9322.5
$ws.Range("M132").Value = -4071.4121
$ws.Range("N132").Value = -14382.5

$ws.Range("H138").Value = 95000
$ws.Range("J138").Value = 95000
$ws.Range("L138").Value = 95000
$ws.Range("N138").Value = -105280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2899.2144
$ws.Range("I7").Value = 3007.3333
$ws.Range("K7").Value = 3007.3333
$ws.Range("M7").Value = -2895.3333

$ws.Range("H40").Value = 1901.5
$ws.Range("I40").Value = 1932.1538
$ws.Range("K40").Value = 1932.1538
$ws.Range("M40").Value = -1796.1538

$ws.Range("H61").Value = 3294
$ws.Range("I61").Value = 3123.4
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 3123.4
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2921.4
$ws.Range("N61").Value = -5404

$ws.Range("H104").Value = 24992.5
$ws.Range("J104").Value = 24992.5
$ws.Range("L104").Value = 24992.5
$ws.Range("N104").Value = -31980.5

$ws.Range("H113").Value = 3294
$ws.Range("I113").Value = 3123.4
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 3123.4
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -953.4000000000001
$ws.Range("N113").Value = -9340

$ws.Range("H126").Value = 2899.2144
$ws.Range("I126").Value = 3007.3333
$ws.Range("K126").Value = 9021.999899999999
$ws.Range("M126").Value = -6551.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H113").Value = 839.1667
$ws.Range("J113").Value = 1210.5714
$ws.Range("L113").Value = 3631.7142
$ws.Range("N113").Value = -7971.7142

$ws.Range("H132").Value = 1940.4038
$ws.Range("I132").Value = 1521.9231
$ws.Range("K132").Value = 4565.7693
$ws.Range("M132").Value = -2035.7693

$ws.Range("H136").Value = 39061.066
$ws.Range("I136").Value = 91054.164
$ws.Range("J136").Value = 4399
$ws.Range("K136").Value = 273162.492
$ws.Range("L136").Value = 13197
$ws.Range("M136").Value = -270612.492
$ws.Range("N136").Value = -18297
